# Proposed Plan 2 sheet updated
#
# - Insert a new worksheet "Sheet3" between "Sheet1" and "Student".
# - Populate "Sheet3" with a Date/Time/Batch/PRN/Name/Present header row
#   and six attendance rows (kept in the same fill order as the source
#   commit so the shared-string table comes out in the same sequence).
# - Move the cursor/selection on "Sheet1" to B5.
# - Leave "Sheet3" as the active sheet/tab with A1:F1 selected.

$wb = $excel.ActiveWorkbook

# --- Sheet1: just move the selection to B5 -------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("B5").Select() | Out-Null

# --- Insert the new "Sheet3" worksheet right after "Sheet1" --------------
$ws3 = $wb.Worksheets.Add($null, $ws1)
$ws3.Name = "Sheet3"

# Header cells that reuse existing shared strings (Date/Time/Batch).
$ws3.Range("A1").Value = "Date"
$ws3.Range("B1").Value = "Time"
$ws3.Range("C1").Value = "Batch"

# Attendance rows (rows 3-8) filled before the remaining headers so the
# shared-string table ends up with the same ordering as the target file.
$ws3.Range("A3").Value = "| 25-02-2026 | 13:30 | Both  | 2501132001 | Patil Yoksh Laxman      | 1 |"
$ws3.Range("A4").Value = "| 25-02-2026 | 13:30 | Both  | 2501132002 | Uttekar Paarth Hanumant | 1 |"
$ws3.Range("A5").Value = "| 25-02-2026 | 13:30 | Both  | 2501132003 | Gupta Aastha Vijay      | 0 |"
$ws3.Range("A6").Value = "| 25-02-2026 | 14:30 | D1    | 2501132001 | Patil Yoksh Laxman      | 1 |"
$ws3.Range("A7").Value = "| 25-02-2026 | 14:30 | D1    | 2501132002 | Uttekar Paarth Hanumant | 0 |"
$ws3.Range("A8").Value = "| 26-02-2026 | 09:30 | D2    | 2501132004 | Satvik Anand            | 1 |"

# Remaining header cells (PRN / Name / Present).
$ws3.Range("D1").Value = "PRN "
$ws3.Range("E1").Value = "Name"
$ws3.Range("F1").Value = "Present"

# Keep "Sheet3" active/selected with A1:F1 highlighted.
$ws3.Activate()
$ws3.Range("A1:F1").Select() | Out-Null
